# finish gemini with grounding
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)
$lo = $ws.ListObjects.Item(1)

# --- Grow the table to its final size (A1:J4) -------------------------------
$lo.Resize($ws.Range("A1:J4"))

# --- Cells that must end up completely empty (no value, no style) ------------
$ws.Range("D2").Clear()
$ws.Range("H3").Clear()

# D3 carried the old wrapped "Used Prompts" formatting - drop it, it gets a
# plain (non-wrapped) value below.
$ws.Range("D3").Clear()

# --- Header row (row 1) ------------------------------------------------------
$ws.Range("A1").Value = "Version"
$ws.Range("B1").Value = "Changes"
$ws.Range("C1").Value = "Model"
$ws.Range("D1").Value = "context"
$ws.Range("E1").Value = "Answer Format"
$ws.Range("F1").Value = "Field prompts"
$ws.Range("G1").Value = "Query for every"
$ws.Range("H1").Value = "Note"
$ws.Range("I1").Value = "Google Search"
$ws.Range("J1").Value = "Complete?"

# --- Row 2 --------------------------------------------------------------------
$ws.Range("A2").Value = "firm_properties_gemini_with_grounding_v1"
$ws.Range("C2").Value = "gemini-1.5-pro-002"
$ws.Range("E2").Value = "answer_format_1"
$ws.Range("F2").Value = " field_to_query_1"
$ws.Range("G2").Value = "Field"
$ws.Range("I2").Value = "Y"
$ws.Range("J2").Value = "-"

# --- Row 3 --------------------------------------------------------------------
$ws.Range("A3").Value = "firm_properties_gemini_with_grounding_v2"
$ws.Range("B3").Value = "Added new fields"
$ws.Range("C3").Value = "gemini-1.5-pro-002"
$ws.Range("D3").Value = "context_single_answer_v1"
$ws.Range("E3").Value = "answer_format_1"
$ws.Range("F3").Value = " field_to_query_1"
$ws.Range("G3").Value = "Firm"
$ws.Range("I3").Value = "Y"
$ws.Range("J3").Value = 2000

# --- Row 4 --------------------------------------------------------------------
$ws.Range("A4").Value = "firm_properties_gemini_without_grounding_local_dataset_v1"
$ws.Range("D4").Value = "context_single_answer_v1 + context_local_dataset_v1"
$ws.Range("J4").Value = "-"

# --- Wrap text styling (matches style index 1 = wrapText) --------------------
$ws.Range("E2:G2").WrapText = $true
$ws.Range("E3:G3").WrapText = $true
$ws.Range("E4:G4").WrapText = $true

# --- Column widths -------------------------------------------------------------
$ws.Columns.Item(1).ColumnWidth = 47.6666666666667
$ws.Columns.Item(3).ColumnWidth = 15.9986979166667
$ws.Columns.Item(4).ColumnWidth = 44.6666666666667
$ws.Columns.Item(5).ColumnWidth = 67.9986979166667
$ws.Columns.Item(6).ColumnWidth = 67.9986979166667
$ws.Columns.Item(7).ColumnWidth = 67.9986979166667

# --- View: selection + scrolled position --------------------------------------
$ws.Range("J5").Select()

# --- Page setup (orientation portrait) -----------------------------------------
$ws.PageSetup.Orientation = 1
